$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.030.47"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.847.31"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "704.97"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.92"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "3.846.37"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.33"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.79"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "4.499.17"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "3.913.79"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "71.091.03"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.80"
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.65"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.36"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000148"
$ws.Range("E25").Value = "  +2.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.63"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.11"
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.20"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.48"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.179"
$ws.Range("E34").Value = "  -5.49%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").Value = "3.804.45"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.37"
$ws.Range("E39").Value = "  +6.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.03"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  +5.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("B43").Value = "FLOKI"
$ws.Range("C43").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000326"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.53"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.59"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "419.48"
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.40"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.61"
$ws.Range("E51").Value = "  +0.33%  "
